$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=5424; B=45969.95833333334},
    @{Row=3;  A=5455; B=45969.96875},
    @{Row=4;  A=5430; B=45969.97916666666},
    @{Row=5;  A=5337; B=45969.98958333334},
    @{Row=6;  A=5216; B=45970},
    @{Row=7;  A=5126; B=45970.01041666666},
    @{Row=8;  A=5129; B=45970.02083333334},
    @{Row=9;  A=5037; B=45970.03125},
    @{Row=10; A=5048; B=45970.04166666666},
    @{Row=11; A=5060; B=45970.05208333334},
    @{Row=12; A=5022; B=45970.0625},
    @{Row=13; A=5029; B=45970.07291666666},
    @{Row=14; A=4966; B=45970.08333333334},
    @{Row=15; A=4994; B=45970.09375},
    @{Row=16; A=4913; B=45970.10416666666},
    @{Row=17; A=4884; B=45970.11458333334},
    @{Row=18; A=4913; B=45970.125},
    @{Row=19; A=4897; B=45970.13541666666},
    @{Row=20; A=4904; B=45970.14583333334},
    @{Row=21; A=4886; B=45970.15625},
    @{Row=22; A=4956; B=45970.16666666666},
    @{Row=23; A=4906; B=45970.17708333334},
    @{Row=24; A=4948; B=45970.1875},
    @{Row=25; A=4986; B=45970.19791666666},
    @{Row=26; A=5046; B=45970.20833333334},
    @{Row=27; A=5055; B=45970.21875},
    @{Row=28; A=5116; B=45970.22916666666},
    @{Row=29; A=5102; B=45970.23958333334},
    @{Row=30; A=5093; B=45970.25},
    @{Row=31; A=5143; B=45970.26041666666},
    @{Row=32; A=5207; B=45970.27083333334},
    @{Row=33; A=5258; B=45970.28125},
    @{Row=34; A=5269; B=45970.29166666666},
    @{Row=35; A=5295; B=45970.30208333334}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
}
